# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price (col D) and Volume(1h) (col E) are stored as text in this sheet, so
# for prices that look like plain numbers we briefly format the cell as Text
# before writing the value (otherwise Excel's .Value setter would coerce a
# numeric-looking string into a real number) and then restore the default
# "Normal" style so no visible formatting changes stick around.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.818.49'
$ws.Range('E2').Value = '  -0.09%  '

$ws.Range('D3').Value = '2.274.90'
$ws.Range('E3').Value = '  -0.12%  '

$ws.Range('E4').Value = '  -0.30%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '249.48'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.76%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.643'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.00%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '76.90'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +7.22%  '

$ws.Range('E8').Value = '  -0.13%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.654'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.76%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '40.13'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.80%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0973'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.04%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.31'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.14%  '

$ws.Range('E13').Value = '  +0.42%  '

$ws.Range('D14').Value = '2.614.62'
$ws.Range('E14').Value = '  -0.24%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.99'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.63%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.865'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.40%  '

$ws.Range('D17').Value = '2.276.60'
$ws.Range('E17').Value = '  +0.24%  '

$ws.Range('D18').Value = '42.726.21'
$ws.Range('E18').Value = '  -0.27%  '

$ws.Range('E19').Value = '  -2.27%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.21'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.53%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '72.01'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.96%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '233.71'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.99%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.14'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.25%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.84'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -4.80%  '

$ws.Range('E25').Value = '  -0.10%  '

$ws.Range('E26').Value = '  -0.99%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.35'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -3.79%  '

$ws.Range('E28').Value = '  +1.63%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '167.64'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.24%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '20.57'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.95%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.41'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.42%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0854'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +6.92%  '

$ws.Range('E33').Value = '  -3.58%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '30.43'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.38%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.127'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.77%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.56'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.09%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.72'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.99%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0303'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -3.40%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '13.73'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +3.19%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.26'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.76%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.85'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.27%  '

$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '112.12'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +17.57%  '

$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.208'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.59%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '60.99'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.73%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.85'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.52%  '

$ws.Range('E46').Value = '  -1.30%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.00'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.30%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '4.50'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -10.09%  '

$ws.Range('E49').Value = '  -2.85%  '

$ws.Range('E50').Value = '  -2.41%  '

$ws.Range('E51').Value = '  -0.08%  '
